$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B9").Value = 3464198.2
$ws.Range("C9").Value = 540091.99
$ws.Range("D9").Value = 4004290.19
$ws.Range("E9").Value = 13.48783340799783
$ws.Range("F9").Value = 86.51216659200216
$ws.Range("G9").Value = -47.80280218931496
$ws.Range("H9").Value = -37.44136885998463
$ws.Range("I9").Value = 34676
$ws.Range("J9").Value = 1468
$ws.Range("K9").Value = 36144
$ws.Range("L9").Value = 24955
$ws.Range("M9").Value = 160.4604363854939
$ws.Range("N9").Value = 9.549608135125954
